# Scheduled market-data refresh: update cached price/profit figures
# (columns H:N) for the affected Leve rows across each crafting-job sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1047.175
$ws.Range("I129").Value = 274.25
$ws.Range("J129").Value = 1133.0555
$ws.Range("K129").Value = 822.75
$ws.Range("L129").Value = 3399.1665
$ws.Range("M129").Value = 4177.25
$ws.Range("N129").Value = -13399.1665
$ws.Range("H138").Value = 5284.231
$ws.Range("I138").Value = 1293.3704
$ws.Range("J138").Value = 8119.8423
$ws.Range("K138").Value = 3880.1112
$ws.Range("L138").Value = 24359.5269
$ws.Range("M138").Value = 1259.8888
$ws.Range("N138").Value = -34639.5269

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19711.848
$ws.Range("I32").Value = 21548.264
$ws.Range("J32").Value = 12224.923
$ws.Range("K32").Value = 21548.264
$ws.Range("L32").Value = 12224.923
$ws.Range("M32").Value = -21261.264
$ws.Range("N32").Value = -12798.923
$ws.Range("H63").Value = 2751.25
$ws.Range("I63").Value = 2502.5
$ws.Range("K63").Value = 2502.5
$ws.Range("M63").Value = -1816.5
$ws.Range("H66").Value = 2751.25
$ws.Range("I66").Value = 2502.5
$ws.Range("K66").Value = 12512.5
$ws.Range("M66").Value = -9080.5
$ws.Range("H88").Value = 6832.154
$ws.Range("I88").Value = 10336.333
$ws.Range("J88").Value = 3828.5715
$ws.Range("K88").Value = 10336.333
$ws.Range("L88").Value = 3828.5715
$ws.Range("M88").Value = -9930.333000000001
$ws.Range("N88").Value = -4640.5715
$ws.Range("H91").Value = 6832.154
$ws.Range("I91").Value = 10336.333
$ws.Range("J91").Value = 3828.5715
$ws.Range("K91").Value = 10336.333
$ws.Range("L91").Value = 3828.5715
$ws.Range("M91").Value = -8932.333000000001
$ws.Range("N91").Value = -6636.5715
$ws.Range("H110").Value = 1406.619
$ws.Range("I110").Value = 1373.9412
$ws.Range("J110").Value = 1545.5
$ws.Range("K110").Value = 1373.9412
$ws.Range("L110").Value = 1545.5
$ws.Range("M110").Value = 671.0588
$ws.Range("N110").Value = -5635.5
$ws.Range("H134").Value = 65422
$ws.Range("J134").Value = 65422
$ws.Range("L134").Value = 65422
$ws.Range("N134").Value = -75562

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1103.25
$ws.Range("I20").Value = 1118
$ws.Range("J20").Value = 1000
$ws.Range("K20").Value = 1118
$ws.Range("L20").Value = 1000
$ws.Range("M20").Value = -871
$ws.Range("N20").Value = -1494
$ws.Range("H22").Value = 340.25
$ws.Range("I22").Value = 330
$ws.Range("J22").Value = 350.5
$ws.Range("K22").Value = 330
$ws.Range("L22").Value = 350.5
$ws.Range("M22").Value = -157
$ws.Range("N22").Value = -696.5
$ws.Range("H86").Value = 1845.5
$ws.Range("I86").Value = 1851.5122
$ws.Range("J86").Value = 1599
$ws.Range("K86").Value = 1851.5122
$ws.Range("L86").Value = 1599
$ws.Range("M86").Value = -728.5121999999999
$ws.Range("N86").Value = -3845
$ws.Range("H89").Value = 1845.5
$ws.Range("I89").Value = 1851.5122
$ws.Range("J89").Value = 1599
$ws.Range("K89").Value = 9257.561
$ws.Range("L89").Value = 7995
$ws.Range("M89").Value = -3641.561
$ws.Range("N89").Value = -19227
$ws.Range("H118").Value = 56039.4
$ws.Range("J118").Value = 56039.4
$ws.Range("L118").Value = 56039.4
$ws.Range("N118").Value = -59353.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1846.92
$ws.Range("I99").Value = 1368.2941
$ws.Range("K99").Value = 1368.2941
$ws.Range("M99").Value = 129.7058999999999
$ws.Range("H126").Value = 1846.92
$ws.Range("I126").Value = 1368.2941
$ws.Range("K126").Value = 4104.8823
$ws.Range("M126").Value = -1634.8823
$ws.Range("H134").Value = 2993.111
$ws.Range("I134").Value = 1889.1052
$ws.Range("K134").Value = 5667.3156
$ws.Range("M134").Value = -3132.3156

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 10701.23
$ws.Range("I3").Value = 7236.6665
$ws.Range("J3").Value = 13670.857
$ws.Range("K3").Value = 21709.9995
$ws.Range("L3").Value = 41012.571
$ws.Range("M3").Value = -21597.9995
$ws.Range("N3").Value = -41236.571

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5793.2793
$ws.Range("I70").Value = 5194.7646
$ws.Range("J70").Value = 6184.615
$ws.Range("K70").Value = 5194.7646
$ws.Range("L70").Value = 6184.615
$ws.Range("M70").Value = -4924.7646
$ws.Range("N70").Value = -6724.615
$ws.Range("H73").Value = 5793.2793
$ws.Range("I73").Value = 5194.7646
$ws.Range("J73").Value = 6184.615
$ws.Range("K73").Value = 5194.7646
$ws.Range("L73").Value = 6184.615
$ws.Range("M73").Value = -4258.7646
$ws.Range("N73").Value = -8056.615
$ws.Range("H80").Value = 8219.6
$ws.Range("I80").Value = 13443.429
$ws.Range("J80").Value = 5406.769
$ws.Range("K80").Value = 13443.429
$ws.Range("L80").Value = 5406.769
$ws.Range("M80").Value = -12445.429
$ws.Range("N80").Value = -7402.769
$ws.Range("H83").Value = 8219.6
$ws.Range("I83").Value = 13443.429
$ws.Range("J83").Value = 5406.769
$ws.Range("K83").Value = 67217.145
$ws.Range("L83").Value = 27033.845
$ws.Range("M83").Value = -62225.145
$ws.Range("N83").Value = -37017.845
$ws.Range("H107").Value = 403.22223
$ws.Range("I107").Value = 147.53847
$ws.Range("K107").Value = 147.53847
$ws.Range("M107").Value = 1772.46153
$ws.Range("H113").Value = 1632.6897
$ws.Range("I113").Value = 1512.2222
$ws.Range("J113").Value = 1829.8182
$ws.Range("K113").Value = 1512.2222
$ws.Range("L113").Value = 1829.8182
$ws.Range("M113").Value = 657.7778000000001
$ws.Range("N113").Value = -6169.8182
$ws.Range("H122").Value = 5486.357
$ws.Range("I122").Value = 7138.125
$ws.Range("J122").Value = 3284
$ws.Range("K122").Value = 21414.375
$ws.Range("L122").Value = 9852
$ws.Range("M122").Value = -18964.375
$ws.Range("N122").Value = -14752

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1835.7142
$ws.Range("I68").Value = 1830
$ws.Range("J68").Value = 1850
$ws.Range("K68").Value = 1830
$ws.Range("L68").Value = 1850
$ws.Range("M68").Value = -1081
$ws.Range("N68").Value = -3348
$ws.Range("H71").Value = 1835.7142
$ws.Range("I71").Value = 1830
$ws.Range("J71").Value = 1850
$ws.Range("K71").Value = 9150
$ws.Range("L71").Value = 9250
$ws.Range("M71").Value = -5406
$ws.Range("N71").Value = -16738
$ws.Range("H82").Value = 1961.5834
$ws.Range("J82").Value = 2231.5
$ws.Range("L82").Value = 2231.5
$ws.Range("N82").Value = -2953.5
$ws.Range("H85").Value = 1961.5834
$ws.Range("J85").Value = 2231.5
$ws.Range("L85").Value = 2231.5
$ws.Range("N85").Value = -4727.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 3093.7334
$ws.Range("I107").Value = 745.55554
$ws.Range("J107").Value = 6616
$ws.Range("K107").Value = 2236.66662
$ws.Range("L107").Value = 19848
$ws.Range("M107").Value = -316.66662
$ws.Range("N107").Value = -23688
